$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-formatted columns (B, C, D) keep their values as text, not numbers
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.307.37'
$ws.Range('E2').Value = '  +0.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.874.27'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7097'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '241.80'
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07862'
$ws.Range('E8').Value = '  +2.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3122'
$ws.Range('E9').Value = '  +0.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '25.16'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08385'
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.874.77'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.236'
$ws.Range('E13').Value = '  +0.53%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.7168'
$ws.Range('E14').Value = '  +1.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.21'
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.208'
$ws.Range('E16').Value = '  +4.29%  '
$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.323.18'
$ws.Range('E17').Value = '  +0.23%  '
$ws.Range('B18').NumberFormat = '@'
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').NumberFormat = '@'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008322'
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.63'
$ws.Range('E19').Value = '  -0.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.22'
$ws.Range('E20').Value = '  +0.28%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.118.95'
$ws.Range('E21').Value = '  -1.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.9999'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.774'
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('E24').Value = '  +0.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1593'
$ws.Range('E25').Value = '  -2.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.044'
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.47'
$ws.Range('E27').Value = '  -0.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.54'
$ws.Range('E28').Value = '  +0.25%  '
$ws.Range('E29').Value = '  +0.24%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.413'
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.338'
$ws.Range('E31').Value = '  +0.58%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.220'
$ws.Range('E32').Value = '  -4.86%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05351'
$ws.Range('E33').Value = '  +2.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.945'
$ws.Range('E34').Value = '  +1.29%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.176'
$ws.Range('E35').Value = '  +0.47%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7471'
$ws.Range('E36').Value = '  -0.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.695'
$ws.Range('E37').Value = '  +0.57%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.288.63'
$ws.Range('E38').Value = '  +11.82%  '
$ws.Range('E39').Value = '  +1.38%  '
$ws.Range('E40').Value = '  +0.95%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.566'
$ws.Range('E41').Value = '  +3.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.8954'
$ws.Range('E42').Value = '  +1.28%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '110.71'
$ws.Range('E43').Value = '  +5.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '72.93'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000132'
$ws.Range('E45').Value = '  +9.57%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.9999'
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.021.23'
$ws.Range('E47').Value = '  -0.22%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.799'
$ws.Range('E48').Value = '  +0.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5195'
$ws.Range('E49').Value = '  +0.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.442'
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4352'
$ws.Range('E51').Value = '  +1.36%  '
